# Auto-generated: update cryptos price/volume figures per commit diff.
# Values are forced to Text (leading apostrophe) then style is reset to
# "Normal" so the stored cell keeps default formatting (no stray numFmt),
# matching the source workbook where these columns are plain text cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'62.792.03"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +8.78%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'3.470.28"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +5.93%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.22%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'414.48"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +3.88%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'124.07"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +13.68%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'3.461.45"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +5.76%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'  +1.93%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'  +0.14%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.681"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +9.43%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.128"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +33.01%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'41.25"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +4.68%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("E13").Value = "'  +0.54%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'4.016.97"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +6.10%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'8.59"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +3.65%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'19.95"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +4.87%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'3.479.33"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +6.32%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'62.768.55"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +9.37%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "'  -0.23%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'10.85"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -2.00%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").Value = "'  +28.53%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = "'  -0.05%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'316.74"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +6.53%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'82.12"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +10.43%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'12.91"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -0.10%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'3.17"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -0.45%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'30.81"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +9.55%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'7.79"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +5.31%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'7.83"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -0.99%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value = "'  +3.32%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'4.31"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -1.99%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'0.118"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +5.03%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'2.66"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +24.18%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'11.55"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +2.54%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'42.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +4.04%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = "'  +0.15%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'0.0492"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -2.35%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'52.20"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +0.59%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'3.51"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +0.61%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value = "'  -0.01%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = "'  -2.78%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'2.00"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +6.26%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "'  +2.98%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'136.05"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -1.75%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.284"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -0.06%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'16.88"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +0.26%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'  -0.50%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = "'  +2.05%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'21.91"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -1.67%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'2.205.85"
$ws.Range("D50").Style = "Normal"
$ws.Range("E51").Value = "'  +0.47%  "
$ws.Range("E51").Style = "Normal"
